{"js": "const replacements = [\n    [\"84\u00f78=10, 4\", \"92\u00f74=23, 0\"],\n    [\"76\u00f74=19, 0\", \"77\u00f75=15, 2\"],\n    [\"34\u00f75=6, 4\", \"72\u00f72=36, 0\"],\n    [\"71\u00f77=10, 1\", \"53\u00f77=7, 4\"],\n    [\"78\u00f76=13, 0\", \"61\u00f72=30, 1\"],\n    [\"55\u00f76=9, 1\", \"82\u00f74=20, 2\"],\n    [\"44\u00f75=8, 4\", \"78\u00f77=11, 1\"],\n    [\"35\u00f72=17, 1\", \"18\u00f75=3, 3\"],\n    [\"30\u00f72=15, 0\", \"85\u00f77=12, 1\"],\n    [\"94\u00f72=47, 0\", \"15\u00f79=1, 6\"],\n    [\"15\u00f75=3, 0\", \"77\u00f79=8, 5\"],\n    [\"37\u00f76=6, 1\", \"67\u00f72=33, 1\"],\n    [\"50\u00f75=10, 0\", \"60\u00f74=15, 0\"],\n    [\"67\u00f78=8, 3\", \"67\u00f72=33, 1\"],\n    [\"28\u00f79=3, 1\", \"78\u00f72=39, 0\"],\n    [\"18\u00f72=9, 0\", \"54\u00f72=27, 0\"],\n    [\"96\u00f75=19, 1\", \"53\u00f72=26, 1\"],\n    [\"73\u00f76=12, 1\", \"67\u00f75=13, 2\"],\n    [\"60\u00f72=30, 0\", \"22\u00f76=3, 4\"],\n    [\"80\u00f79=8, 8\", \"74\u00f78=9, 2\"],\n    [\"87\u00f72=43, 1\", \"67\u00f73=22, 1\"],\n    [\"70\u00f74=17, 2\", \"57\u00f77=8, 1\"],\n    [\"57\u00f73=19, 0\", \"16\u00f74=4, 0\"],\n    [\"21\u00f76=3, 3\", \"97\u00f75=19, 2\"],\n    [\"52\u00f76=8, 4\", \"20\u00f74=5, 0\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('84\u00f78=10, 4', '92\u00f74=23, 0'),\n    @('76\u00f74=19, 0', '77\u00f75=15, 2'),\n    @('34\u00f75=6, 4', '72\u00f72=36, 0'),\n    @('71\u00f77=10, 1', '53\u00f77=7, 4'),\n    @('78\u00f76=13, 0', '61\u00f72=30, 1'),\n    @('55\u00f76=9, 1', '82\u00f74=20, 2'),\n    @('44\u00f75=8, 4', '78\u00f77=11, 1'),\n    @('35\u00f72=17, 1', '18\u00f75=3, 3'),\n    @('30\u00f72=15, 0', '85\u00f77=12, 1'),\n    @('94\u00f72=47, 0', '15\u00f79=1, 6'),\n    @('15\u00f75=3, 0', '77\u00f79=8, 5'),\n    @('37\u00f76=6, 1', '67\u00f72=33, 1'),\n    @('50\u00f75=10, 0', '60\u00f74=15, 0'),\n    @('67\u00f78=8, 3', '67\u00f72=33, 1'),\n    @('28\u00f79=3, 1', '78\u00f72=39, 0'),\n    @('18\u00f72=9, 0', '54\u00f72=27, 0'),\n    @('96\u00f75=19, 1', '53\u00f72=26, 1'),\n    @('73\u00f76=12, 1', '67\u00f75=13, 2'),\n    @('60\u00f72=30, 0', '22\u00f76=3, 4'),\n    @('80\u00f79=8, 8', '74\u00f78=9, 2'),\n    @('87\u00f72=43, 1', '67\u00f73=22, 1'),\n    @('70\u00f74=17, 2', '57\u00f77=8, 1'),\n    @('57\u00f73=19, 0', '16\u00f74=4, 0'),\n    @('21\u00f76=3, 3', '97\u00f75=19, 2'),\n    @('52\u00f76=8, 4', '20\u00f74=5, 0')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
